$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")
# ALC row 4 (hunk 0)
$ws_ALC.Range("H4").Value = 207.875
$ws_ALC.Range("I4").Value = 166.14285
$ws_ALC.Range("K4").Value = 166.14285
$ws_ALC.Range("M4").Value = -52.14285000000001

# ALC row 32 (hunk 1)
$ws_ALC.Range("H32").Value = 19444.445
$ws_ALC.Range("J32").Value = 21000
$ws_ALC.Range("L32").Value = 21000
$ws_ALC.Range("N32").Value = -21652

# ALC row 33 (hunk 2)
$ws_ALC.Range("H33").Value = 419
$ws_ALC.Range("I33").Value = 414.9091
$ws_ALC.Range("K33").Value = 414.9091
$ws_ALC.Range("M33").Value = -185.9091

# ALC row 96 (hunk 3)
$ws_ALC.Range("H96").Value = 2080.2
$ws_ALC.Range("I96").Value = 2581.5
$ws_ALC.Range("K96").Value = 7744.5
$ws_ALC.Range("M96").Value = -6371.5

# ARM row 2 (hunk 4)
$ws_ARM.Range("H2").Value = 1164.5358
$ws_ARM.Range("I2").Value = 937.2105
$ws_ARM.Range("J2").Value = 1644.4445
$ws_ARM.Range("K2").Value = 937.2105
$ws_ARM.Range("L2").Value = 1644.4445
$ws_ARM.Range("M2").Value = -824.2105
$ws_ARM.Range("N2").Value = -1870.4445

# ARM row 63 (hunk 5)
$ws_ARM.Range("H63").Value = 0
$ws_ARM.Range("I63").Value = 0
$ws_ARM.Range("K63").Value = 0
$ws_ARM.Range("M63").ClearContents()

# ARM row 66 (hunk 6)
$ws_ARM.Range("H66").Value = 0
$ws_ARM.Range("I66").Value = 0
$ws_ARM.Range("K66").Value = 0
$ws_ARM.Range("M66").ClearContents()

# ARM row 116 (hunk 7)
$ws_ARM.Range("H116").Value = 1164.5358
$ws_ARM.Range("I116").Value = 937.2105
$ws_ARM.Range("J116").Value = 1644.4445
$ws_ARM.Range("K116").Value = 937.2105
$ws_ARM.Range("L116").Value = 1644.4445
$ws_ARM.Range("M116").Value = 1356.7895
$ws_ARM.Range("N116").Value = -6232.4445

# ARM row 132 (hunk 8)
$ws_ARM.Range("H132").Value = 7088.613
$ws_ARM.Range("I132").Value = 6165.6313
$ws_ARM.Range("J132").Value = 8550
$ws_ARM.Range("K132").Value = 18496.8939
$ws_ARM.Range("L132").Value = 25650
$ws_ARM.Range("M132").Value = -15966.8939
$ws_ARM.Range("N132").Value = -30710

# BSM row 3 (hunk 9)
$ws_BSM.Range("H3").Value = 1164.5358
$ws_BSM.Range("I3").Value = 937.2105
$ws_BSM.Range("J3").Value = 1644.4445
$ws_BSM.Range("K3").Value = 937.2105
$ws_BSM.Range("L3").Value = 1644.4445
$ws_BSM.Range("M3").Value = -823.2105
$ws_BSM.Range("N3").Value = -1872.4445

# BSM row 76 (hunk 10)
$ws_BSM.Range("H76").Value = 8362.666999999999
$ws_BSM.Range("J76").Value = 8362.666999999999
$ws_BSM.Range("L76").Value = 8362.666999999999
$ws_BSM.Range("N76").Value = -8992.666999999999

# BSM row 79 (hunk 11)
$ws_BSM.Range("H79").Value = 8362.666999999999
$ws_BSM.Range("J79").Value = 8362.666999999999
$ws_BSM.Range("L79").Value = 8362.666999999999
$ws_BSM.Range("N79").Value = -10546.667

# BSM row 134 (hunk 12)
$ws_BSM.Range("H134").Value = 1472.8182
$ws_BSM.Range("I134").Value = 1472.8182
$ws_BSM.Range("K134").Value = 4418.4546
$ws_BSM.Range("M134").Value = -1883.4546

# CRP row 10 (hunk 13)
$ws_CRP.Range("H10").Value = 439
$ws_CRP.Range("I10").Value = 439
$ws_CRP.Range("J10").Value = 0
$ws_CRP.Range("K10").Value = 439
$ws_CRP.Range("L10").Value = 0
$ws_CRP.Range("M10").Value = -300
$ws_CRP.Range("N10").ClearContents()

# CRP row 31 (hunk 14)
$ws_CRP.Range("H31").Value = 6647
$ws_CRP.Range("I31").Value = 9947
$ws_CRP.Range("K31").Value = 9947
$ws_CRP.Range("M31").Value = -9652

# CRP row 34 (hunk 15)
$ws_CRP.Range("H34").Value = 6647
$ws_CRP.Range("I34").Value = 9947
$ws_CRP.Range("K34").Value = 9947
$ws_CRP.Range("M34").Value = -9745

# CRP row 36 (hunk 16)
$ws_CRP.Range("H36").Value = 0
$ws_CRP.Range("I36").Value = 0
$ws_CRP.Range("K36").Value = 0
$ws_CRP.Range("M36").ClearContents()

# CRP row 40 (hunk 17)
$ws_CRP.Range("H40").Value = 0
$ws_CRP.Range("I40").Value = 0
$ws_CRP.Range("K40").Value = 0
$ws_CRP.Range("M40").ClearContents()

# CRP row 42 (hunk 18)
$ws_CRP.Range("H42").Value = 19500
$ws_CRP.Range("I42").Value = 0
$ws_CRP.Range("J42").Value = 19500
$ws_CRP.Range("K42").Value = 0
$ws_CRP.Range("L42").Value = 19500
$ws_CRP.Range("M42").ClearContents()
$ws_CRP.Range("N42").Value = -20686

# CRP row 44 (hunk 19)
$ws_CRP.Range("H44").Value = 0
$ws_CRP.Range("I44").Value = 0
$ws_CRP.Range("K44").Value = 0
$ws_CRP.Range("M44").ClearContents()

# CRP row 62 (hunk 20)
$ws_CRP.Range("H62").Value = 5000
$ws_CRP.Range("I62").Value = 0
$ws_CRP.Range("J62").Value = 5000
$ws_CRP.Range("K62").Value = 0
$ws_CRP.Range("L62").Value = 5000
$ws_CRP.Range("M62").ClearContents()
$ws_CRP.Range("N62").Value = -6248

# CRP row 65 (hunk 21)
$ws_CRP.Range("H65").Value = 5000
$ws_CRP.Range("I65").Value = 0
$ws_CRP.Range("J65").Value = 5000
$ws_CRP.Range("K65").Value = 0
$ws_CRP.Range("L65").Value = 25000
$ws_CRP.Range("M65").ClearContents()
$ws_CRP.Range("N65").Value = -31240

# CRP row 86 (hunk 22)
$ws_CRP.Range("H86").Value = 14123.125
$ws_CRP.Range("I86").Value = 9665
$ws_CRP.Range("K86").Value = 9665
$ws_CRP.Range("M86").Value = -8542

# CRP row 89 (hunk 23)
$ws_CRP.Range("H89").Value = 14123.125
$ws_CRP.Range("I89").Value = 9665
$ws_CRP.Range("K89").Value = 48325
$ws_CRP.Range("M89").Value = -42709

# CUL row 2 (hunk 24)
$ws_CUL.Range("H2").Value = 794.1429000000001
$ws_CUL.Range("J2").Value = 1365.25
$ws_CUL.Range("L2").Value = 8191.5
$ws_CUL.Range("N2").Value = -8417.5

# CUL row 41 (hunk 25)
$ws_CUL.Range("H41").Value = 488.5
$ws_CUL.Range("I41").Value = 230
$ws_CUL.Range("K41").Value = 690
$ws_CUL.Range("M41").Value = -352

# CUL row 62 (hunk 26)
$ws_CUL.Range("H62").Value = 8849.5
$ws_CUL.Range("I62").Value = 1100
$ws_CUL.Range("J62").Value = 10399.4
$ws_CUL.Range("K62").Value = 3300
$ws_CUL.Range("L62").Value = 31198.2
$ws_CUL.Range("M62").Value = -2614
$ws_CUL.Range("N62").Value = -32570.2

# CUL row 65 (hunk 27)
$ws_CUL.Range("H65").Value = 8849.5
$ws_CUL.Range("I65").Value = 1100
$ws_CUL.Range("J65").Value = 10399.4
$ws_CUL.Range("K65").Value = 9900
$ws_CUL.Range("L65").Value = 93594.59999999999
$ws_CUL.Range("M65").Value = -6468
$ws_CUL.Range("N65").Value = -100458.6

# CUL row 102 (hunk 28)
$ws_CUL.Range("H102").Value = 750
$ws_CUL.Range("J102").Value = 750
$ws_CUL.Range("L102").Value = 2250
$ws_CUL.Range("N102").Value = -7118

# CUL row 108 (hunk 29)
$ws_CUL.Range("H108").Value = 1845
$ws_CUL.Range("I108").Value = 1845
$ws_CUL.Range("K108").Value = 5535
$ws_CUL.Range("M108").Value = -2655

# GSM row 3 (hunk 30)
$ws_GSM.Range("H3").Value = 30000000
$ws_GSM.Range("I3").Value = 30000000
$ws_GSM.Range("K3").Value = 30000000
$ws_GSM.Range("M3").Value = -29999884

# GSM row 7 (hunk 31)
$ws_GSM.Range("H7").Value = 5002
$ws_GSM.Range("J7").Value = 0
$ws_GSM.Range("L7").Value = 0
$ws_GSM.Range("N7").ClearContents()

# GSM row 8 (hunk 32)
$ws_GSM.Range("H8").Value = 5002
$ws_GSM.Range("J8").Value = 0
$ws_GSM.Range("L8").Value = 0
$ws_GSM.Range("N8").ClearContents()

# GSM row 10 (hunk 33)
$ws_GSM.Range("H10").Value = 507497.5
$ws_GSM.Range("J10").Value = 14995
$ws_GSM.Range("L10").Value = 14995
$ws_GSM.Range("N10").Value = -15333

# GSM row 70 (hunk 34)
$ws_GSM.Range("H70").Value = 1000
$ws_GSM.Range("I70").Value = 1000
$ws_GSM.Range("K70").Value = 1000
$ws_GSM.Range("M70").Value = -730

# GSM row 73 (hunk 35)
$ws_GSM.Range("H73").Value = 1000
$ws_GSM.Range("I73").Value = 1000
$ws_GSM.Range("K73").Value = 1000
$ws_GSM.Range("M73").Value = -64

# GSM row 102 (hunk 36)
$ws_GSM.Range("H102").Value = 1816.6072
$ws_GSM.Range("I102").Value = 1406.875
$ws_GSM.Range("J102").Value = 4275
$ws_GSM.Range("K102").Value = 1406.875
$ws_GSM.Range("L102").Value = 4275
$ws_GSM.Range("M102").Value = 215.125
$ws_GSM.Range("N102").Value = -7519

# GSM row 104 (hunk 37)
$ws_GSM.Range("H104").Value = 52379.145
$ws_GSM.Range("J104").Value = 52379.145
$ws_GSM.Range("L104").Value = 52379.145
$ws_GSM.Range("N104").Value = -59367.145

# LTW row 7 (hunk 38)
$ws_LTW.Range("H7").Value = 3522.3333
$ws_LTW.Range("I7").Value = 3525.125
$ws_LTW.Range("J7").Value = 3500
$ws_LTW.Range("K7").Value = 3525.125
$ws_LTW.Range("L7").Value = 3500
$ws_LTW.Range("M7").Value = -3413.125
$ws_LTW.Range("N7").Value = -3724

# LTW row 22 (hunk 39)
$ws_LTW.Range("H22").Value = 1480.1666
$ws_LTW.Range("I22").Value = 1470.75
$ws_LTW.Range("J22").Value = 1499
$ws_LTW.Range("K22").Value = 1470.75
$ws_LTW.Range("L22").Value = 1499
$ws_LTW.Range("M22").Value = -1175.75
$ws_LTW.Range("N22").Value = -2089

# LTW row 27 (hunk 40)
$ws_LTW.Range("H27").Value = 1480.1666
$ws_LTW.Range("I27").Value = 1470.75
$ws_LTW.Range("J27").Value = 1499
$ws_LTW.Range("K27").Value = 1470.75
$ws_LTW.Range("L27").Value = 1499
$ws_LTW.Range("M27").Value = -1363.75
$ws_LTW.Range("N27").Value = -1713

# LTW row 46 (hunk 41)
$ws_LTW.Range("H46").Value = 6636.3335
$ws_LTW.Range("I46").Value = 4800
$ws_LTW.Range("J46").Value = 7554.5
$ws_LTW.Range("K46").Value = 4800
$ws_LTW.Range("L46").Value = 7554.5
$ws_LTW.Range("M46").Value = -4612
$ws_LTW.Range("N46").Value = -7930.5

# LTW row 55 (hunk 42)
$ws_LTW.Range("H55").Value = 1133.9231
$ws_LTW.Range("I55").Value = 388.6
$ws_LTW.Range("J55").Value = 1599.75
$ws_LTW.Range("K55").Value = 388.6
$ws_LTW.Range("L55").Value = 1599.75
$ws_LTW.Range("M55").Value = -215.6
$ws_LTW.Range("N55").Value = -1945.75

# LTW row 126 (hunk 43)
$ws_LTW.Range("H126").Value = 3522.3333
$ws_LTW.Range("I126").Value = 3525.125
$ws_LTW.Range("J126").Value = 3500
$ws_LTW.Range("K126").Value = 10575.375
$ws_LTW.Range("L126").Value = 10500
$ws_LTW.Range("M126").Value = -8105.375
$ws_LTW.Range("N126").Value = -15440

# LTW row 132 (hunk 44)
$ws_LTW.Range("H132").Value = 3239.9666
$ws_LTW.Range("I132").Value = 3116.625
$ws_LTW.Range("K132").Value = 9349.875
$ws_LTW.Range("M132").Value = -6819.875

# WVR row 3 (hunk 45)
$ws_WVR.Range("H3").Value = 25000000
$ws_WVR.Range("I3").Value = 25000000
$ws_WVR.Range("J3").Value = 0
$ws_WVR.Range("K3").Value = 25000000
$ws_WVR.Range("L3").Value = 0
$ws_WVR.Range("M3").Value = -24999886
$ws_WVR.Range("N3").ClearContents()

# WVR row 4 (hunk 46)
$ws_WVR.Range("H4").Value = 8700
$ws_WVR.Range("J4").Value = 8700
$ws_WVR.Range("L4").Value = 8700
$ws_WVR.Range("N4").Value = -8926

# WVR row 41 (hunk 47)
$ws_WVR.Range("H41").Value = 16998.5
$ws_WVR.Range("I41").Value = 16998.5
$ws_WVR.Range("K41").Value = 16998.5
$ws_WVR.Range("M41").Value = -16608.5

# WVR row 75 (hunk 48)
$ws_WVR.Range("H75").Value = 57000
$ws_WVR.Range("J75").Value = 57000
$ws_WVR.Range("L75").Value = 57000
$ws_WVR.Range("N75").Value = -58872

# WVR row 78 (hunk 49)
$ws_WVR.Range("H78").Value = 57000
$ws_WVR.Range("J78").Value = 57000
$ws_WVR.Range("L78").Value = 171000
$ws_WVR.Range("N78").Value = -180360

# WVR row 96 (hunk 50)
$ws_WVR.Range("H96").Value = 575.5
$ws_WVR.Range("I96").Value = 367.33334
$ws_WVR.Range("J96").Value = 1200
$ws_WVR.Range("K96").Value = 367.33334
$ws_WVR.Range("L96").Value = 1200
$ws_WVR.Range("M96").Value = 1005.66666
$ws_WVR.Range("N96").Value = -3946

# WVR row 136 (hunk 51)
$ws_WVR.Range("H136").Value = 38152.863
$ws_WVR.Range("I136").Value = 36970.777
$ws_WVR.Range("J136").Value = 43472.25
$ws_WVR.Range("K136").Value = 110912.331
$ws_WVR.Range("L136").Value = 130416.75
$ws_WVR.Range("M136").Value = -108362.331
$ws_WVR.Range("N136").Value = -135516.75
